$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IP")

# Insert a new column E (pushes old "connection" column E -> F)
$ws.Range("E1").EntireColumn.Insert()
$ws.Cells.Item(1,5).Value = "ssh port"
$ws.Columns.Item(5).ColumnWidth = 11.804421768707465

# Set ssh port values for existing router rows (2-21)
$ws.Cells.Item(2,5).Value = 22101
$ws.Cells.Item(5,5).Value = 22102
$ws.Cells.Item(10,5).Value = 22103
$ws.Cells.Item(14,5).Value = 22104
$ws.Cells.Item(18,5).Value = 22105

# Insert 8 new rows before row 22 for dmzdns/dmzmail/dmzweb/proxy
$ws.Range("A22:A29").EntireRow.Insert()

# Insert 6 new rows before (old row28=entwks101, now shifted to row36) for mail/web/file
$ws.Range("A36:A41").EntireRow.Insert()

# Fill in rows 22-43 with full target content
$ws.Cells.Item(22,1).Value = "dmzdns"
$ws.Cells.Item(22,2).Value = "em0"
$ws.Cells.Item(22,3).Value = "OPT1"
$ws.Cells.Item(22,4).Value = "nat"
$ws.Cells.Item(22,5).Value = 22301
$ws.Cells.Item(23,1).Value = "dmzdns"
$ws.Cells.Item(23,2).Value = "em1"
$ws.Cells.Item(23,3).Value = "DMZ"
$ws.Cells.Item(23,4).Value = "10.10.3.2"
$ws.Cells.Item(24,1).Value = "dmzmail"
$ws.Cells.Item(24,2).Value = "em0"
$ws.Cells.Item(24,3).Value = "OPT1"
$ws.Cells.Item(24,4).Value = "nat"
$ws.Cells.Item(24,5).Value = 22302
$ws.Cells.Item(25,1).Value = "dmzmail"
$ws.Cells.Item(25,2).Value = "em1"
$ws.Cells.Item(25,3).Value = "DMZ"
$ws.Cells.Item(25,4).Value = "10.10.3.3"
$ws.Cells.Item(26,1).Value = "dmzweb"
$ws.Cells.Item(26,2).Value = "em0"
$ws.Cells.Item(26,3).Value = "OPT1"
$ws.Cells.Item(26,4).Value = "nat"
$ws.Cells.Item(26,5).Value = 22304
$ws.Cells.Item(27,1).Value = "dmzweb"
$ws.Cells.Item(27,2).Value = "em1"
$ws.Cells.Item(27,3).Value = "DMZ"
$ws.Cells.Item(27,4).Value = "10.10.3.4"
$ws.Cells.Item(28,1).Value = "proxy"
$ws.Cells.Item(28,2).Value = "em0"
$ws.Cells.Item(28,3).Value = "OPT1"
$ws.Cells.Item(28,4).Value = "nat"
$ws.Cells.Item(28,5).Value = 22303
$ws.Cells.Item(29,1).Value = "proxy"
$ws.Cells.Item(29,2).Value = "em1"
$ws.Cells.Item(29,3).Value = "DMZ"
$ws.Cells.Item(29,4).Value = "10.10.3.5"
$ws.Cells.Item(30,1).Value = "nagios"
$ws.Cells.Item(30,2).Value = "em0"
$ws.Cells.Item(30,3).Value = "OPT1"
$ws.Cells.Item(30,4).Value = "nat"
$ws.Cells.Item(30,5).Value = 22401
$ws.Cells.Item(30,6).Value = "http://127.0.0.1:18033/nagios/"
$ws.Cells.Item(31,1).Value = "nagios"
$ws.Cells.Item(31,2).Value = "em1"
$ws.Cells.Item(31,3).Value = "SOC"
$ws.Cells.Item(31,4).Value = "10.10.4.3"
$ws.Cells.Item(32,1).Value = "socws101"
$ws.Cells.Item(32,2).Value = "em0"
$ws.Cells.Item(32,3).Value = "OPT1"
$ws.Cells.Item(32,4).Value = "nat"
$ws.Cells.Item(32,5).Value = 22402
$ws.Cells.Item(33,1).Value = "socws101"
$ws.Cells.Item(33,2).Value = "em1"
$ws.Cells.Item(33,3).Value = "SOC"
$ws.Cells.Item(33,4).Value = "10.10.4.101"
$ws.Cells.Item(34,1).Value = "DC1"
$ws.Cells.Item(34,2).Value = "em0"
$ws.Cells.Item(34,3).Value = "OPT1"
$ws.Cells.Item(34,4).Value = "nat"
$ws.Cells.Item(34,5).Value = 22801
$ws.Cells.Item(35,1).Value = "DC1"
$ws.Cells.Item(35,2).Value = "em1"
$ws.Cells.Item(35,3).Value = "ENTSVR"
$ws.Cells.Item(35,4).Value = "10.10.8.11"
$ws.Cells.Item(36,1).Value = "mail"
$ws.Cells.Item(36,2).Value = "em0"
$ws.Cells.Item(36,3).Value = "OPT1"
$ws.Cells.Item(36,4).Value = "nat"
$ws.Cells.Item(36,5).Value = 22903
$ws.Cells.Item(37,1).Value = "mail"
$ws.Cells.Item(37,2).Value = "em1"
$ws.Cells.Item(37,3).Value = "ENTSVR"
$ws.Cells.Item(37,4).Value = "10.10.8.13"
$ws.Cells.Item(38,1).Value = "web"
$ws.Cells.Item(38,2).Value = "em0"
$ws.Cells.Item(38,3).Value = "OPT1"
$ws.Cells.Item(38,4).Value = "nat"
$ws.Cells.Item(38,5).Value = 22904
$ws.Cells.Item(39,1).Value = "web"
$ws.Cells.Item(39,2).Value = "em1"
$ws.Cells.Item(39,3).Value = "ENTSVR"
$ws.Cells.Item(39,4).Value = "10.10.8.14"
$ws.Cells.Item(40,1).Value = "file"
$ws.Cells.Item(40,2).Value = "em0"
$ws.Cells.Item(40,3).Value = "OPT1"
$ws.Cells.Item(40,4).Value = "nat"
$ws.Cells.Item(40,5).Value = 22905
$ws.Cells.Item(41,1).Value = "file"
$ws.Cells.Item(41,2).Value = "em1"
$ws.Cells.Item(41,3).Value = "ENTSVR"
$ws.Cells.Item(41,4).Value = "10.10.8.15"
$ws.Cells.Item(42,1).Value = "entwks101"
$ws.Cells.Item(42,2).Value = "em0"
$ws.Cells.Item(42,3).Value = "OPT1"
$ws.Cells.Item(42,4).Value = "nat"
$ws.Cells.Item(42,5).Value = 22901
$ws.Cells.Item(43,1).Value = "entwks101"
$ws.Cells.Item(43,2).Value = "em1"
$ws.Cells.Item(43,3).Value = "ENTWKS"
$ws.Cells.Item(43,4).Value = "10.10.9.101"

# Update selection to match target
$ws.Range("E22").Select() | Out-Null
